$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2615.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2615.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 7845.75
$ws.Range("N17").Value = -8181.75
$ws.Range("H92").Value = 50524.45
$ws.Range("I92").Value = 62870.75
$ws.Range("J92").Value = 1139.25
$ws.Range("K92").Value = 62870.75
$ws.Range("L92").Value = 1139.25
$ws.Range("M92").Value = -61622.75
$ws.Range("N92").Value = -3635.25
$ws.Range("H98").Value = 1589.1428
$ws.Range("I98").Value = 1287.2727
$ws.Range("J98").Value = 2696
$ws.Range("K98").Value = 1287.2727
$ws.Range("L98").Value = 2696
$ws.Range("M98").Value = 210.7273
$ws.Range("N98").Value = -5692
$ws.Range("H106").Value = 7999.7334
$ws.Range("I106").Value = 6076.75
$ws.Range("J106").Value = 15691.667
$ws.Range("K106").Value = 6076.75
$ws.Range("L106").Value = 15691.667
$ws.Range("M106").Value = -5445.75
$ws.Range("N106").Value = -16953.667
$ws.Range("H107").Value = 2393
$ws.Range("I107").Value = 2265.375
$ws.Range("J107").Value = 2733.3333
$ws.Range("K107").Value = 2265.375
$ws.Range("L107").Value = 2733.3333
$ws.Range("M107").Value = -345.375
$ws.Range("N107").Value = -6573.3333
$ws.Range("H111").Value = 1217.5714
$ws.Range("I111").Value = 1217.5714
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3652.7142
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -585.7142000000003
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 1589.1428
$ws.Range("I122").Value = 1287.2727
$ws.Range("J122").Value = 2696
$ws.Range("K122").Value = 3861.8181
$ws.Range("L122").Value = 8088
$ws.Range("M122").Value = -1411.8181
$ws.Range("N122").Value = -12988
$ws.Range("H127").Value = 2499.75
$ws.Range("I127").Value = 2499.75
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 7499.25
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -2539.25
$ws.Range("H138").Value = 4362.838
$ws.Range("I138").Value = 2987.842
$ws.Range("J138").Value = 5814.222
$ws.Range("K138").Value = 8963.526
$ws.Range("L138").Value = 17442.666
$ws.Range("M138").Value = -3823.526
$ws.Range("N138").Value = -27722.666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2527.182
$ws.Range("I61").Value = 1314.4286
$ws.Range("J61").Value = 4649.5
$ws.Range("K61").Value = 1314.4286
$ws.Range("L61").Value = 4649.5
$ws.Range("M61").Value = -1102.4286
$ws.Range("N61").Value = -5073.5
$ws.Range("H74").Value = 40019.15
$ws.Range("I74").Value = 41404.5
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 41404.5
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -40530.5
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 40019.15
$ws.Range("I77").Value = 41404.5
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 207022.5
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -202654.5
$ws.Range("N77").Value = -28736
$ws.Range("H122").Value = 1920.6428
$ws.Range("I122").Value = 1606.8462
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 4820.5386
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -2370.5386
$ws.Range("N122").Value = -22900
$ws.Range("H136").Value = 2527.182
$ws.Range("I136").Value = 1314.4286
$ws.Range("J136").Value = 4649.5
$ws.Range("K136").Value = 3943.2858
$ws.Range("L136").Value = 13948.5
$ws.Range("M136").Value = -1393.2858
$ws.Range("N136").Value = -19048.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2844.375
$ws.Range("I86").Value = 2648.4546
$ws.Range("J86").Value = 4999.5
$ws.Range("K86").Value = 2648.4546
$ws.Range("L86").Value = 4999.5
$ws.Range("M86").Value = -1525.4546
$ws.Range("N86").Value = -7245.5
$ws.Range("H89").Value = 2844.375
$ws.Range("I89").Value = 2648.4546
$ws.Range("J89").Value = 4999.5
$ws.Range("K89").Value = 13242.273
$ws.Range("L89").Value = 24997.5
$ws.Range("M89").Value = -7626.273000000001
$ws.Range("N89").Value = -36229.5
$ws.Range("H134").Value = 1795.6346
$ws.Range("I134").Value = 1477.9783
$ws.Range("J134").Value = 4231
$ws.Range("K134").Value = 4433.9349
$ws.Range("L134").Value = 12693
$ws.Range("M134").Value = -1898.9349
$ws.Range("N134").Value = -17763

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2200.5356
$ws.Range("I31").Value = 2068
$ws.Range("J31").Value = 2995.75
$ws.Range("K31").Value = 2068
$ws.Range("L31").Value = 2995.75
$ws.Range("M31").Value = -1773
$ws.Range("N31").Value = -3585.75
$ws.Range("H34").Value = 2200.5356
$ws.Range("I34").Value = 2068
$ws.Range("J34").Value = 2995.75
$ws.Range("K34").Value = 2068
$ws.Range("L34").Value = 2995.75
$ws.Range("M34").Value = -1866
$ws.Range("N34").Value = -3399.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 987.0833
$ws.Range("I34").Value = 987.0833
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2961.2499
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2877.2499
$ws.Range("N34").ClearContents()
$ws.Range("H55").Value = 3529.3333
$ws.Range("I55").Value = 3529.3333
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 10587.9999
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -10410.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 37833
$ws.Range("I44").Value = 38500
$ws.Range("J44").Value = 37499.5
$ws.Range("K44").Value = 38500
$ws.Range("L44").Value = 37499.5
$ws.Range("M44").Value = -37904
$ws.Range("N44").Value = -38691.5
$ws.Range("H47").Value = 39249.5
$ws.Range("I47").Value = 38500
$ws.Range("J47").Value = 39999
$ws.Range("K47").Value = 38500
$ws.Range("L47").Value = 39999
$ws.Range("M47").Value = -37932
$ws.Range("N47").Value = -41135
$ws.Range("H122").Value = 3096.4
$ws.Range("I122").Value = 2827.8333
$ws.Range("J122").Value = 3499.25
$ws.Range("K122").Value = 8483.499899999999
$ws.Range("L122").Value = 10497.75
$ws.Range("M122").Value = -6033.499899999999
$ws.Range("N122").Value = -15397.75
$ws.Range("H126").Value = 7051
$ws.Range("I126").Value = 6722.2
$ws.Range("J126").Value = 7462
$ws.Range("K126").Value = 20166.6
$ws.Range("L126").Value = 22386
$ws.Range("M126").Value = -17696.6
$ws.Range("N126").Value = -27326
$ws.Range("H132").Value = 26084.953
$ws.Range("I132").Value = 35430.7
$ws.Range("J132").Value = 2720.5833
$ws.Range("K132").Value = 106292.1
$ws.Range("L132").Value = 8161.749899999999
$ws.Range("M132").Value = -103762.1
$ws.Range("N132").Value = -13221.7499
$ws.Range("H137").Value = 73000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 73000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 73000
$ws.Range("N137").Value = -83200

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2989.3
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 3982.1667
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3982.1667
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -4572.1667
$ws.Range("H27").Value = 2989.3
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 3982.1667
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3982.1667
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -4196.1667
$ws.Range("H115").Value = 79999
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 79999
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 79999
$ws.Range("N115").Value = -82349
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H136").Value = 3522.7646
$ws.Range("I136").Value = 2590.4167
$ws.Range("J136").Value = 5760.4
$ws.Range("K136").Value = 7771.250100000001
$ws.Range("L136").Value = 17281.2
$ws.Range("M136").Value = -5221.250100000001
$ws.Range("N136").Value = -22381.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 53419
$ws.Range("I45").Value = 7997
$ws.Range("J45").Value = 64774.5
$ws.Range("K45").Value = 7997
$ws.Range("L45").Value = 64774.5
$ws.Range("M45").Value = -7506
$ws.Range("N45").Value = -65756.5

Write-Output "applied all changes"